# "Generate Report for Handback"
#
# The localization-status report is regenerated: the cf35edde-... file now
# has come back from handback ("Handed back: in sync with en-US") and moves
# into the first data row (row 2) on every sheet, picking up its Target /
# Handback file + datetime columns. The c93166b9-... file (still only
# "Ready for handoff") slides down to row 3. Row 1 headers are untouched.

$wb = $excel.ActiveWorkbook

# Hyperlink font used throughout the sheet for "looks like a link" cells
# (matches the workbook's existing HyperLink cell style: underlined, blue).
$HL_COLOR = 15570276   # OLE BGR for #6495ED

function Style-AsLink($rng) {
    $rng.Font.Underline = $true
    $rng.Font.Color = $HL_COLOR
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# The two existing hyperlinks (A2 -> c93166b9, A3 -> cf35edde) keep pointing
# at the same targets (same r:id in the diff) - only the displayed text
# swaps along with the cell text, so row 2 now reads as the cf35edde file
# and row 3 as the c93166b9 file.
$ovLinks = @()
foreach ($hl in $ws.Hyperlinks) { $ovLinks += $hl }
foreach ($hl in $ovLinks) {
    $addr = $hl.Range.Address($false, $false)
    if ($addr -eq "A2") {
        $hl.TextToDisplay = "cf35edde-5631-431f-bfd9-e86b1269f077.md"
    } elseif ($addr -eq "A3") {
        $hl.TextToDisplay = "c93166b9-b152-4ed2-9188-d7e4c736fc46.md"
    }
}

$ws.Range("A2").Value2 = "cf35edde-5631-431f-bfd9-e86b1269f077.md"
$ws.Range("B2").Value2 = "Handed back: in sync with en-US"
$ws.Range("C2").Value2 = "Handed back: in sync with en-US"
$ws.Range("D2").Value2 = "2016-31-12 04:31:30"

$ws.Range("A3").Value2 = "c93166b9-b152-4ed2-9188-d7e4c736fc46.md"
$ws.Range("B3").Value2 = "Ready for handoff"
$ws.Range("C3").Value2 = "Ready for handoff"
$ws.Range("D3").Value2 = "2016-31-12 04:31:15"

# ---------------------------------------------------------------------
# Per-locale sheets ("zh-cn" and "de-de"): same row swap, plus the
# cf35edde row gains "Latest Target File" (F) / "Latest Handback File" (G)
# and a real "Latest Handback DateTime" (H).
# ---------------------------------------------------------------------
$locales = @(
    @{
        Sheet   = "zh-cn"
        Ext     = "zh-cn.xlf"
        CfXlf   = "cf35edde-5631-431f-bfd9-e86b1269f077.4cc94648b413eb21f705f6e44db0fd7c87e3509d.zh-cn.xlf"
        C93Xlf  = "c93166b9-b152-4ed2-9188-d7e4c736fc46.608a58dbbce996c93cec27acc58bd782e7ac473d.zh-cn.xlf"
        CfXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6e28b54347ab4f1f2bf98b533a07f1b7e5eeb2bc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/cf35edde-5631-431f-bfd9-e86b1269f077.4cc94648b413eb21f705f6e44db0fd7c87e3509d.zh-cn.xlf"
        C93XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea9db38bbe98d22b0ea1e7c79ef69c99352f5001/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c93166b9-b152-4ed2-9188-d7e4c736fc46.608a58dbbce996c93cec27acc58bd782e7ac473d.zh-cn.xlf"
        CfHandoffDt  = "2016-03-12 04:31:27"
        CfHandbackDt = "2016-03-12 04:31:40"
        C93HandoffDt = "2016-03-12 04:31:12"
    },
    @{
        Sheet   = "de-de"
        Ext     = "de-de.xlf"
        CfXlf   = "cf35edde-5631-431f-bfd9-e86b1269f077.4cc94648b413eb21f705f6e44db0fd7c87e3509d.de-de.xlf"
        C93Xlf  = "c93166b9-b152-4ed2-9188-d7e4c736fc46.608a58dbbce996c93cec27acc58bd782e7ac473d.de-de.xlf"
        CfXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ed438f9b328fbdf643f14ac6142d010a63d8b87d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/cf35edde-5631-431f-bfd9-e86b1269f077.4cc94648b413eb21f705f6e44db0fd7c87e3509d.de-de.xlf"
        C93XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dda4701c6815516b21b38d4316b2f066d5dafbfa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c93166b9-b152-4ed2-9188-d7e4c736fc46.608a58dbbce996c93cec27acc58bd782e7ac473d.de-de.xlf"
        CfHandoffDt  = "2016-03-12 04:31:30"
        CfHandbackDt = "2016-03-12 04:31:46"
        C93HandoffDt = "2016-03-12 04:31:15"
    }
)

$CfMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/1e1ed23483cd12e60bd2879d4e747503cc0c626d/e2e/cf35edde-5631-431f-bfd9-e86b1269f077.md"
$C93MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/e953b6a709a03d41f8760a231a4d2a8ffd407a45/e2e/c93166b9-b152-4ed2-9188-d7e4c736fc46.md"
$CfMdName  = "cf35edde-5631-431f-bfd9-e86b1269f077.md"
$C93MdName = "c93166b9-b152-4ed2-9188-d7e4c736fc46.md"

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    # Re-point the 6 existing hyperlinks (A2,B2,D2,A3,B3,D3) at their new
    # targets - row 2 now belongs to cf35edde, row 3 to c93166b9.
    $links = @()
    foreach ($hl in $ws.Hyperlinks) { $links += $hl }
    foreach ($hl in $links) {
        $addr = $hl.Range.Address($false, $false)
        switch ($addr) {
            "A2" { $hl.Address = $CfMdUrl;       $hl.TextToDisplay = $CfMdName }
            "B2" { $hl.Address = $CfMdUrl;       $hl.TextToDisplay = ".md" }
            "D2" { $hl.Address = $loc.CfXlfUrl;  $hl.TextToDisplay = $loc.CfXlf }
            "A3" { $hl.Address = $C93MdUrl;      $hl.TextToDisplay = $C93MdName }
            "B3" { $hl.Address = $C93MdUrl;      $hl.TextToDisplay = ".md" }
            "D3" { $hl.Address = $loc.C93XlfUrl; $hl.TextToDisplay = $loc.C93Xlf }
        }
    }

    # New hyperlinks for the cf35edde row's Target File (F2) / Handback
    # File (G2) columns.
    $ws.Hyperlinks.Add($ws.Range("F2"), $CfMdUrl, "", "", $CfMdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $loc.CfXlfUrl, "", "", $loc.CfXlf) | Out-Null
    Style-AsLink($ws.Range("F2"))
    Style-AsLink($ws.Range("G2"))

    # Row 2 <- cf35edde (now "handed back")
    $ws.Range("A2").Value2 = $CfMdName
    $ws.Range("B2").Value2 = ".md"
    $ws.Range("C2").Value2 = "Handed back: in sync with en-US"
    $ws.Range("D2").Value2 = $loc.CfXlf
    $ws.Range("E2").Value2 = $loc.CfHandoffDt
    $ws.Range("F2").Value2 = $CfMdName
    $ws.Range("G2").Value2 = $loc.CfXlf
    $ws.Range("H2").Value2 = $loc.CfHandbackDt
    $ws.Range("I2").Value2 = "Include"

    # Row 3 <- c93166b9 (still just "ready for handoff")
    $ws.Range("A3").Value2 = $C93MdName
    $ws.Range("B3").Value2 = ".md"
    $ws.Range("C3").Value2 = "Ready for handoff"
    $ws.Range("D3").Value2 = $loc.C93Xlf
    $ws.Range("E3").Value2 = $loc.C93HandoffDt
    $ws.Range("H3").Value2 = "0001-01-01 00:00:00"
    $ws.Range("I3").Value2 = "Include"
}
